# carry forward hourly run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Increment the TestCases run number (B2) and the Instance tag (D2)
$ws.Range("B2").Value = "55"
$ws.Range("D2").Value = "Automation2"

# Move the active selection to D2, matching the edited sheet view
$ws.Range("D2").Select()
